# The sheet holds a header row (row 1) plus two "helper" rows (rows 2 and 3)
# that contained example/placeholder values ("Campo livre", a sample process
# number mask, "dd/mm/aaaa", "*obrigatório", the allowed action-type list,
# etc.). Those example values are being removed, leaving the rows empty but
# still formatted - this matches the commit's cleanup of the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) also had a handful of trailing cells (V1:Z1) whose style
# differed only by a missing "bold" application versus the rest of the
# header - bring them in line with the rest of the header row.
$ws.Range("V1:Z1").Font.Bold = $true

# Clear the sample/placeholder content from the two helper rows while
# keeping their cell formatting intact.
$ws.Range("A2:U2").ClearContents() | Out-Null
$ws.Range("A3:Z3").ClearContents() | Out-Null

# Reflect the resulting selection (the range that was cleared).
$ws.Range("A2:XFD3").Select() | Out-Null
